$d = $word.ActiveDocument

# Vertical-tab character == <w:br/> manual line break when it shows up inside
# Range text / Find-replace text in this object model.
$BR = [char]11

# NB: we deliberately do *not* pass $replaceText straight to Find.Execute's
# own "Replace" argument - that path runs the replacement string through
# Word's AutoFormat/AutoCorrect ("smart quotes") logic, silently turning a
# plain apostrophe (') into a curly one (U+2019). Instead we use Find only
# to *locate* the range, then assign Range.Text directly, which substitutes
# the text verbatim while still inheriting the first run's formatting and
# still honouring embedded vertical-tab (manual line break) characters.
function Replace-Text($findText, $replaceText) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Find.Execute could not locate: " + $findText)
    }
    $rng.Text = $replaceText
    return $found
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Text "Exploring the Wonders of the Quantum Realm" "Embracing Diversity in a Globalized World: Navigating Challenges and Celebrating Unity"

# ---------------------------------------------------------------------------
# Author name (collapses the five runs making up "Dr. Albert J. Robertson"
# down to the single new run "Hannah Davies")
# ---------------------------------------------------------------------------
Replace-Text "Dr. Albert J. Robertson" "Hannah Davies"

# ---------------------------------------------------------------------------
# Author email
# ---------------------------------------------------------------------------
Replace-Text "quantumstudies@scientificdiscovery" "hannah.davies@eduworld"

# ---------------------------------------------------------------------------
# Body paragraph - first block (before the first manual line break)
# ---------------------------------------------------------------------------
$find1 = "In the vast expanse of the cosmos, beyond the boundaries of our visible reality, lies a realm where particles dance in an intricate choreography, governed by the inexplicable laws of quantum mechanics."
$find1 = $find1 + " This hidden world, teeming with enigmatic phenomena and possibilities, beckons curious minds to unravel its mysteries and unlock the secrets that hold the key to our understanding of the universe."
$replace1 = "In a world where borders are increasingly blurred and cultures intertwine, embracing diversity has become a cornerstone of global harmony."
$replace1 = $replace1 + " As citizens of a globalized society, it is imperative that we understand, appreciate, and celebrate the rich tapestry of cultures, beliefs, and perspectives that make our world a vibrant and dynamic place."
$replace1 = $replace1 + " In this essay, we will delve into the multifaceted nature of diversity, explore the challenges it presents, and highlight the immense benefits that stem from fostering inclusivity and understanding."
Replace-Text $find1 $replace1

# ---------------------------------------------------------------------------
# Body paragraph - second block (between 1st and 2nd manual line breaks)
# ---------------------------------------------------------------------------
$find2 = "As we delve into the intricacies of the quantum realm, we transcend the familiar world of classical physics and encounter an arena where particles behave in perplexing and counterintuitive ways."
$find2 = $find2 + " Quantum entanglement, a baffling phenomenon, allows particles to share information instantaneously across vast distances, defying the limitations of space and time."
$find2 = $find2 + " Enter the realm of superposition, where particles exist in multiple states simultaneously, blurring the boundaries between reality and possibility."
$replace2 = "The beauty of diversity lies in its multifaceted nature."
$replace2 = $replace2 + " It encompasses differences in race, ethnicity, gender, sexual orientation, religion, language, and socioeconomic status, among others."
$replace2 = $replace2 + " Each individual brings a unique set of experiences, perspectives, and talents to the table, creating a rich and dynamic social fabric."
$replace2 = $replace2 + " Diversity challenges us to step outside our comfort zones, question our assumptions, and engage with those who hold different beliefs and values."
$replace2 = $replace2 + " By doing so, we broaden our horizons, foster empathy, and gain a deeper understanding of the human experience."
Replace-Text $find2 $replace2

# ---------------------------------------------------------------------------
# Body paragraph - third block (after the 2nd manual line break); this is
# also where a large amount of brand-new material (two "Introduction
# Continued:" sections) gets appended, all still inside the same paragraph,
# separated by pairs of manual line breaks just like the existing ones.
# ---------------------------------------------------------------------------
$find3 = "The quantum realm unveils a universe brimming with uncertainty and indeterminacy, challenging our notions of causality and predictability."
$find3 = $find3 + " The uncertainty principle, a fundamental tenet of this quantum world, dictates that certain properties of particles, such as their position and momentum, cannot be simultaneously known with absolute certainty."
$find3 = $find3 + " This principle introduces an element of inherent fuzziness into our understanding of the underlying fabric of reality."

$replace3 = "However, embracing diversity is not without its challenges."
$replace3 = $replace3 + " Misunderstandings, prejudice, and discrimination can arise when people from different backgrounds interact."
$replace3 = $replace3 + " These challenges can be daunting, but they also present opportunities for growth and learning."
$replace3 = $replace3 + " By promoting open dialogue, encouraging inclusivity, and challenging stereotypes, we can create a society where diversity is celebrated and everyone feels valued and respected."
$replace3 = $replace3 + $BR + $BR + "Introduction Continued:" + $BR
$replace3 = $replace3 + $BR + "The immense benefits that stem from fostering inclusivity and understanding are undeniable."
$replace3 = $replace3 + " A diverse society is a more vibrant, innovative, and resilient one."
$replace3 = $replace3 + " When people from different backgrounds come together, they bring fresh ideas, perspectives, and solutions to complex problems."
$replace3 = $replace3 + " This cross-pollination of ideas leads to "
$replace3 = $replace3 + "groundbreaking discoveries, transformative technologies, and creative expressions that benefit all of humanity."
$replace3 = $replace3 + " Furthermore, a diverse society is a more tolerant and peaceful one."
$replace3 = $replace3 + " When people understand and respect each other's differences, they are less likely to resort to violence or conflict."
$replace3 = $replace3 + " Instead, they are more likely to work together to build a better future for themselves and for generations to come."
$replace3 = $replace3 + $BR + $BR + "Introduction Continued:" + $BR
$replace3 = $replace3 + $BR + "The journey towards a truly diverse and inclusive society is an ongoing one, but it is a journey worth taking."
$replace3 = $replace3 + " By challenging our biases, educating ourselves about different cultures, and actively promoting inclusivity, we can create a world where everyone feels valued, respected, and empowered."
$replace3 = $replace3 + " A world where diversity is not just tolerated but celebrated, and where the unique contributions of each individual are recognized and appreciated."
$replace3 = $replace3 + " In this world, the challenges of diversity will be outweighed by the boundless opportunities it presents, leading to a more harmonious, equitable, and prosperous global community."

Replace-Text $find3 $replace3

# ---------------------------------------------------------------------------
# Summary heading paragraph stays "Summary" (unchanged) - only the body below
# it changes.
# ---------------------------------------------------------------------------
$find4 = "The quantum realm, a place of mystery and wonder, where particles engage in a ceaseless dance governed by the enigmatic laws of quantum mechanics."
$find4 = $find4 + " From quantum entanglement to the principles of superposition and uncertainty, this realm challenges our understanding of causality and predictability, inviting us to contemplate the profound implications of its existence."
$find4 = $find4 + " Delving into The Wonders of the Quantum Realm signifies a journey into the heart of creation itself, where we glimpse the interconnectedness of all things and the boundless possibilities that lie beyond the veil of our ordinary perception."

$replace4 = "In conclusion, embracing diversity in a globalized world presents both challenges and immense benefits."
$replace4 = $replace4 + " By fostering inclusivity, understanding, and respect, we can create a society where everyone feels valued and empowered."
$replace4 = $replace4 + " The challenges of diversity, such as misunderstandings and prejudice, can be overcome through open dialogue, education, and a commitment to building a more just and equitable world."
$replace4 = $replace4 + " The beauty of diversity lies in its multifaceted nature, and it is through celebrating our differences and coming together as a global community that we can unlock the full potential of our shared humanity."

Replace-Text $find4 $replace4

# ---------------------------------------------------------------------------
# Trailing empty paragraph added at the very end of the document body.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
